$wb = $excel.ActiveWorkbook

# New row of data (price date + same-day price, repeated from the previous row)
$newDate = "2025-03-07"

# Map of worksheet name -> price value to append as row 6 (column B)
$updates = @{
    "N-Dense" = "40"
    "N-Type"  = "43"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    $rngA = $ws.Range("A6")
    $rngB = $ws.Range("B6")

    # Force the new cells to Text so the date-like / number-like strings are
    # stored verbatim as text (matching the existing "Date" column and the
    # text-stored prices used on the other sheets), not auto-converted to a
    # real date/number by Excel's input parsing.
    $rngA.NumberFormat = "@"
    $rngB.NumberFormat = "@"

    $rngA.Value = $newDate
    $rngB.Value = $updates[$sheetName]

    # Reset back to the default "Normal" cell style so no stray number
    # format / style index is left behind on the new cells.
    $ws.Range("A6:B6").Style = "Normal"
}
